$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Range("A$row")
    $current = $cell.Value()
    $cell.Value = $current - 1
}
